$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 53 totals (revised figures)
$ws.Range("H53").Value = 34096
$ws.Range("J53").Value = 34096
$ws.Range("N53").Value = 63109
$ws.Range("P53").Value = 61985

# Add new row 54 (01-04-2021 period)
# The date-like label must stay plain text (matches the other period
# labels in column A), so force a text entry via the leading apostrophe
# and then drop back to the Normal style so no extra number format is
# left attached to the cell.
$dateLabel = [string]([char]39) + "01-04-2021"
$ws.Range("A54").Value = $dateLabel
$ws.Range("A54").Style = "Normal"

$ws.Range("B54").Value = 8996
$ws.Range("C54").Value = 767
$ws.Range("D54").Value = 8229
$ws.Range("E54").Value = 3045
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 3045
$ws.Range("H54").Value = 34584
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 34584
$ws.Range("K54").Value = 20120
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 20120
$ws.Range("N54").Value = 66745
$ws.Range("O54").Value = 767
$ws.Range("P54").Value = 65977
